$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.317.38"
$ws.Range("E2").Value = "  -3.55%  "

# Row 3
$ws.Range("D3").Value = "3.717.17"
$ws.Range("E3").Value = "  -1.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.38%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.21"
$ws.Range("E5").Value = "  -4.00%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.35"
$ws.Range("E6").Value = "  -4.77%  "

# Row 7
$ws.Range("D7").Value = "3.716.48"
$ws.Range("E7").Value = "  -1.33%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  -2.68%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.157"
$ws.Range("E10").Value = "  -5.72%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.12"
$ws.Range("E11").Value = "  -7.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -5.26%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.22"
$ws.Range("E13").Value = "  -7.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000239"
$ws.Range("E14").Value = "  -5.87%  "

# Row 15
$ws.Range("D15").Value = "4.336.96"
$ws.Range("E15").Value = "  -0.94%  "

# Row 16
$ws.Range("D16").Value = "3.698.19"
$ws.Range("E16").Value = "  -1.32%  "

# Row 17
$ws.Range("D17").Value = "67.160.21"
$ws.Range("E17").Value = "  -3.88%  "

# Row 18
$ws.Range("E18").Value = "  -5.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.08"
$ws.Range("E19").Value = "  -5.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.98"
$ws.Range("E20").Value = "  -2.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "482.75"
$ws.Range("E21").Value = "  -3.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.92"
$ws.Range("E22").Value = "  -3.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.707"
$ws.Range("E23").Value = "  -2.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.88"
$ws.Range("E24").Value = "  -3.73%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  -12.40%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000136"
$ws.Range("E26").Value = "  +1.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.10"
$ws.Range("E27").Value = "  -6.76%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.13"
$ws.Range("E28").Value = "  -11.90%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.89"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.37"
$ws.Range("E31").Value = "  -4.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.06"
$ws.Range("E32").Value = "  +5.07%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.63"
$ws.Range("E33").Value = "  -7.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  -5.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.59%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.987"
$ws.Range("E36").Value = "  -5.82%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.133"
$ws.Range("E37").Value = "  -3.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("E38").Value = "  -8.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.320"
$ws.Range("E39").Value = "  -9.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "443.78"
$ws.Range("E40").Value = "  -1.97%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.60"
$ws.Range("E41").Value = "  -2.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.97"
$ws.Range("E42").Value = "  -4.98%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.79"
$ws.Range("E43").Value = "  -9.38%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.21"
$ws.Range("E44").Value = "  -4.48%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.13"
$ws.Range("E45").Value = "  -10.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "140.51"
$ws.Range("E46").Value = "  +1.19%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "2.775.90"
$ws.Range("E47").Value = "  -6.06%  "

# Row 48
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0344"
$ws.Range("E49").Value = "  -4.62%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.57"
$ws.Range("E50").Value = "  -6.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.85"
$ws.Range("E51").Value = "  +7.18%  "
